$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.47"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'0.52%"
$ws.Range("E2").ClearFormats()
$ws.Range("E3").Value = "'3.68%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'5.085"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'0.86%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.05591"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'-0.20%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'6.470"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'0.8133"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-0.21%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.8449"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'1.28%"
$ws.Range("E8").ClearFormats()
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1333"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'-0.39%"
$ws.Range("E9").ClearFormats()
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.06977"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'0.38%"
$ws.Range("E10").ClearFormats()
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.02820"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'-0.36%"
$ws.Range("E11").ClearFormats()
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09382"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'-0.19%"
$ws.Range("E12").ClearFormats()
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001515"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-0.15%"
$ws.Range("E13").ClearFormats()
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0006004"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-93.81%"
$ws.Range("E14").ClearFormats()
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006203"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-0.50%"
$ws.Range("E15").ClearFormats()
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.608"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'3.09%"
$ws.Range("E16").ClearFormats()
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.021"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'0.10%"
$ws.Range("E17").ClearFormats()
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.055"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'-1.73%"
$ws.Range("E18").ClearFormats()
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3112"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'-2.31%"
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'0.03178"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'-1.74%"
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'-1.35%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'3.762"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'0.51%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.04654"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'-0.43%"
$ws.Range("E23").ClearFormats()
$ws.Range("D25").Value = "'0.001244"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'0.09%"
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'0.004563"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'6.51%"
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'0.00009601"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'-1.02%"
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'0.0001939"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'-0.04%"
$ws.Range("E28").ClearFormats()
$ws.Range("D40").Value = "'0.03661"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'-0.01%"
$ws.Range("E40").ClearFormats()
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006161"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-0.63%"
$ws.Range("E41").ClearFormats()
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1053"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'0.40%"
$ws.Range("E42").ClearFormats()
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002558"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'-6.26%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.008045"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-1.82%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.00005390"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'1.74%"
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'0.01%"
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'0.1451"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'-19.39%"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.002405"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'19.28%"
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "'0.01%"
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'0.01%"
$ws.Range("E50").ClearFormats()
